$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.908057570457458
$ws.Range("B1").Value = 3.664223432540894
$ws.Range("C1").Value = 2.633979320526123
$ws.Range("D1").Value = 0.9479213356971741
$ws.Range("E1").Value = 0.6215535402297974
